{"js": "// Update the \"Git state\" status paragraph: replace the old pending-push\n// sentence with the new commit-hash line, and add a new paragraph right\n// after it that states the onboarding flow is committed locally.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst targetText = \"- Onboarding flow update is local and pending push.\";\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === targetText) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find target paragraph: \" + targetText);\n}\n\n// Replace the text in place so the run's existing formatting (font,\n// size, etc.) is preserved.\ntarget.insertText(\"- Latest local commit (pending push): d80bdb1\", Word.InsertLocation.replace);\n\n// Insert the new paragraph right after it, with matching formatting.\nconst newPara = target.insertParagraph(\n  \"- Onboarding flow update is committed locally and pending push.\",\n  Word.InsertLocation.after\n);\nnewPara.font.load(\"name,size\");\ntarget.font.load(\"name,size\");\nawait context.sync();\n\nnewPara.font.name = target.font.name;\nnewPara.font.size = target.font.size;\nawait context.sync();\n", "ps1": "# Update requirement status docs for onboarding flow:\n#  - Replace the \"local and pending push\" status line with the\n#    latest local commit hash line.\n#  - Add a new line right after it noting the onboarding flow update\n#    is committed locally and pending push.\n\n$d = $word.ActiveDocument\n\n# Replace the old status sentence with the new commit-hash line,\n# preserving the run formatting of the paragraph it lives in.\n$find = $d.Content.Find\n$find.Execute(\n    \"- Onboarding flow update is local and pending push.\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"- Latest local commit (pending push): d80bdb1\",\n    2\n)\n\n# Insert a new paragraph right after that one (inherits its formatting)\n# and set its text to the new onboarding status line.\n$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)\n$lastPara.Range.InsertParagraphAfter()\n$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)\n$newPara.Range.Text = \"- Onboarding flow update is committed locally and pending push.\"\n"}
